# edit.ps1 - applies the "dele setPriority i to funksjoner" edit to the To do - liste.docx
#
# Summary of changes (see diff):
#  1. "Dele " -> "Fikse " in the "Dele orders_setPriorityDirectionAndReturnIfOrders..." bullet.
#  2. Two new yellow/plain bullets inserted around the
#     "Endre bruk i NOT_MOVING_AT_FLOOR og NOT_MOVING_BETWEEN_FLOORS" bullet:
#        - before it: "Dele i to funksjoner setPriorityDirection og existOrders" (yellow highlight)
#        - after it : "Fikse logikk i setPriorityDirection" (no highlight)
#  3. highlight yellow -> green on "Lage global variabel..." and
#     "Lage funksjon for å endre..." bullets.
#  4. The _GoBack bookmark moves from the end of the
#     "Fikse ett eller annet rundt case AT_FLOOR..." paragraph to the very start of the
#     "Endre bruk i NOT_MOVING_AT_FLOOR..." paragraph (an artifact of where the edits were made).
#  5. lastRenderedPageBreak shifts: removed from the "While" run, added to the
#     "Fikse ett eller annet rundt " run and to the "Ta imot bestilling:" run.

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. "Dele " -> "Fikse "
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "setPriorityDirectionAndReturnIfOrders"
$p = $d.Paragraphs($idx)
$p.Range.Find.Execute("Dele ", $true, $false, $false, $false, $false, $true, 1, $false, "Fikse ", 2)

# ---------------------------------------------------------------------------
# 2. highlight yellow -> green on the two "Lage ..." bullets
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "Lage global variabel"
$d.Paragraphs($idx).Range.Font.HighlightColorIndex = 4  # wdBrightGreen -> "green"

$idx = Find-ParagraphIndex $d "Lage funksjon for"
$d.Paragraphs($idx).Range.Font.HighlightColorIndex = 4  # wdBrightGreen -> "green"

# ---------------------------------------------------------------------------
# 3. Insert the two new bullets around "Endre bruk i NOT_MOVING_AT_FLOOR..."
# ---------------------------------------------------------------------------
$wordXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-WordPkg($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$idx = Find-ParagraphIndex $d "Endre bruk i"
$target = $d.Paragraphs($idx)

# 3a. insert new paragraph BEFORE it: "Dele i to funksjoner setPriorityDirection og existOrders" (yellow)
$rngBefore = $target.Range.Duplicate
$rngBefore.Collapse(1)
$rngBefore.InsertParagraphBefore()

$idx = Find-ParagraphIndex $d "Endre bruk i"
$newParaBeforeIdx = $idx - 1
$bodyXml1 = "<w:p $wordXmlNs><w:pPr><w:pStyle w:val=""Listeavsnitt""/><w:numPr><w:ilvl w:val=""1""/><w:numId w:val=""1""/></w:numPr><w:spacing w:line=""480"" w:lineRule=""auto""/><w:rPr><w:highlight w:val=""yellow""/><w:lang w:val=""nb-NO""/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val=""yellow""/><w:lang w:val=""nb-NO""/></w:rPr><w:t xml:space=""preserve"">Dele i to funksjoner </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:rPr><w:highlight w:val=""yellow""/><w:lang w:val=""nb-NO""/></w:rPr><w:t>setPriorityDirection</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:rPr><w:highlight w:val=""yellow""/><w:lang w:val=""nb-NO""/></w:rPr><w:t xml:space=""preserve""> og </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:rPr><w:highlight w:val=""yellow""/><w:lang w:val=""nb-NO""/></w:rPr><w:t>existOrders</w:t></w:r><w:proofErr w:type=""spellEnd""/></w:p>"
$d.Paragraphs($newParaBeforeIdx).Range.InsertXML((New-WordPkg $bodyXml1))

# 3b. insert new paragraph AFTER it: "Fikse logikk i setPriorityDirection" (plain)
$idx = Find-ParagraphIndex $d "Endre bruk i"
$target = $d.Paragraphs($idx)
$rngAfter = $target.Range.Duplicate
$rngAfter.Collapse(0)
$rngAfter.InsertParagraphAfter()

$idx2 = $idx + 1
$bodyXml2 = "<w:p $wordXmlNs><w:pPr><w:pStyle w:val=""Listeavsnitt""/><w:numPr><w:ilvl w:val=""1""/><w:numId w:val=""1""/></w:numPr><w:spacing w:line=""480"" w:lineRule=""auto""/><w:rPr><w:lang w:val=""nb-NO""/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=""nb-NO""/></w:rPr><w:t xml:space=""preserve"">Fikse logikk i </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:rPr><w:lang w:val=""nb-NO""/></w:rPr><w:t>setPriorityDirection</w:t></w:r><w:proofErr w:type=""spellEnd""/></w:p>"
$d.Paragraphs($idx2).Range.InsertXML((New-WordPkg $bodyXml2))

Write-Host "done stage 3b"

# ---------------------------------------------------------------------------
# 4. Move the _GoBack bookmark
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$idx = Find-ParagraphIndex $d "Endre bruk i"
$target = $d.Paragraphs($idx)
$rngBm = $target.Range.Duplicate
$rngBm.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rngBm)

Write-Host "done stage 4"

# ---------------------------------------------------------------------------
# 5. lastRenderedPageBreak shuffle
# ---------------------------------------------------------------------------

# 5a. add lastRenderedPageBreak to the start of "Fikse ett eller annet rundt " run
#     (this paragraph no longer carries the _GoBack bookmark - it moved away in stage 4)
$idx = Find-ParagraphIndex $d "Fikse ett eller annet rundt"
$p = $d.Paragraphs($idx)
$bodyXml3 = "<w:p $wordXmlNs w:rsidR=""00F07A1B"" w:rsidRPr=""005121B5"" w:rsidRDefault=""00F07A1B"" w:rsidP=""00DD0DA4""><w:pPr><w:pStyle w:val=""Listeavsnitt""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr><w:spacing w:line=""480"" w:lineRule=""auto""/><w:ind w:left=""714"" w:hanging=""357""/><w:rPr><w:lang w:val=""nb-NO""/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=""nb-NO""/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space=""preserve"">Fikse ett eller annet rundt </w:t></w:r><w:r w:rsidR=""00926C8B""><w:rPr><w:lang w:val=""nb-NO""/></w:rPr><w:t xml:space=""preserve"">case </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r w:rsidR=""00926C8B""><w:rPr><w:lang w:val=""nb-NO""/></w:rPr><w:t>AT_FLOOR</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r w:rsidR=""00926C8B""><w:rPr><w:lang w:val=""nb-NO""/></w:rPr><w:t xml:space=""preserve""> i </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r w:rsidR=""00926C8B""><w:rPr><w:lang w:val=""nb-NO""/></w:rPr><w:t>esm.c</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r w:rsidR=""00926C8B""><w:rPr><w:lang w:val=""nb-NO""/></w:rPr><w:t>. Noe gir feilmelding rundt der (manglende parentes eller noe)</w:t></w:r></w:p>"
$p.Range.InsertXML((New-WordPkg $bodyXml3))

Write-Host "done stage 5a"

# 5b. add lastRenderedPageBreak to the "Ta imot bestilling:" run
$idx = Find-ParagraphIndex $d "Ta imot bestilling"
$p = $d.Paragraphs($idx)
$bodyXml4 = "<w:p $wordXmlNs w:rsidR=""001353BE"" w:rsidRPr=""00214518"" w:rsidRDefault=""001353BE"" w:rsidP=""001353BE""><w:pPr><w:rPr><w:highlight w:val=""green""/><w:lang w:val=""nb-NO""/></w:rPr></w:pPr><w:r w:rsidRPr=""00214518""><w:rPr><w:highlight w:val=""green""/><w:lang w:val=""nb-NO""/></w:rPr><w:lastRenderedPageBreak/><w:t>Ta imot bestilling:</w:t></w:r></w:p>"
$p.Range.InsertXML((New-WordPkg $bodyXml4))

Write-Host "done stage 5b"

# 5c. remove lastRenderedPageBreak from the "While" run
$idx = Find-ParagraphIndex $d "While"
$p = $d.Paragraphs($idx)
$bodyXml5 = "<w:p $wordXmlNs w:rsidR=""001353BE"" w:rsidRPr=""00214518"" w:rsidRDefault=""001353BE"" w:rsidP=""001353BE""><w:pPr><w:ind w:left=""720"" w:firstLine=""720""/><w:rPr><w:highlight w:val=""green""/><w:lang w:val=""nb-NO""/></w:rPr></w:pPr><w:proofErr w:type=""spellStart""/><w:r w:rsidRPr=""00214518""><w:rPr><w:highlight w:val=""green""/><w:lang w:val=""nb-NO""/></w:rPr><w:t>While</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r w:rsidRPr=""00214518""><w:rPr><w:highlight w:val=""green""/><w:lang w:val=""nb-NO""/></w:rPr><w:t xml:space=""preserve""> løkke</w:t></w:r></w:p>"
$p.Range.InsertXML((New-WordPkg $bodyXml5))

Write-Host "done stage 5c"

